$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Grade")

# Week 6 overview: add grading comments in column E.
# Set in this order so new shared strings land at the indices the
# original author's save produced (Good, Good choices, Correct, Correctly done).
$ws.Range("E6").Value = "Good"
$ws.Range("E8").Value = "Good choices"
$ws.Range("E7").Value = "Correct"
$ws.Range("E10").Value = "Correctly done"
$ws.Range("E11").Value = "Correct"
$ws.Range("E12").Value = "Correct"

# Widen column E so the new comments are readable.
$ws.Columns("E").ColumnWidth = 21.5

# The rows that now carry a comment grow slightly taller to fit the text.
$ws.Rows("6:6").RowHeight = 17
$ws.Rows("7:7").RowHeight = 17
$ws.Rows("8:8").RowHeight = 17
$ws.Rows("10:10").RowHeight = 17
$ws.Rows("11:11").RowHeight = 17
$ws.Rows("12:12").RowHeight = 17

# Page orientation explicitly set to portrait.
$ws.PageSetup.Orientation = 1

# Leave the selection where the author's last edit (E15) was.
$ws.Range("E15").Select() | Out-Null
